$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'67.395.81"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.21%  "

# Row 3
$ws.Range("D3").Value = "'3.309.52"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.08%  "

# Row 5
$ws.Range("D5").Value = "'185.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.13%  "

# Row 6
$ws.Range("D6").Value = "'577.26"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.09%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("E9").Value = "  -1.04%  "

# Row 10
$ws.Range("E10").Value = "  +0.94%  "

# Row 11
$ws.Range("E11").Value = "  +0.20%  "

# Row 12
$ws.Range("D12").Value = "'3.888.81"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.16%  "

# Row 13
$ws.Range("E13").Value = "  -0.23%  "

# Row 14
$ws.Range("E14").Value = "  -0.20%  "

# Row 15
$ws.Range("D15").Value = "'67.631.27"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.83%  "

# Row 16
$ws.Range("E16").Value = "  -1.03%  "

# Row 17
$ws.Range("D17").Value = "'3.323.27"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.67%  "

# Row 18
$ws.Range("D18").Value = "'443.16"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.88%  "

# Row 19
$ws.Range("D19").Value = "'5.70"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.42%  "

# Row 20
$ws.Range("E20").Value = "  +1.10%  "

# Row 21
$ws.Range("D21").Value = "'7.76"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.49%  "

# Row 22
$ws.Range("D22").Value = "'74.04"
$ws.Range("D22").ClearFormats()

# Row 23
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.05%  "

# Row 24
$ws.Range("D24").Value = "'0.517"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.46%  "

# Row 25
$ws.Range("D25").Value = "'3.456.66"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.16%  "

# Row 26
$ws.Range("E26").Value = "  +0.71%  "

# Row 27
$ws.Range("D27").Value = "'0.188"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.59%  "

# Row 28
$ws.Range("D28").Value = "'9.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.54%  "

# Row 29
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.16%  "

# Row 30
$ws.Range("D30").Value = "'1.97"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.11%  "

# Row 31
$ws.Range("D31").Value = "'22.94"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.80%  "

# Row 32
$ws.Range("D32").Value = "'5.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.94%  "

# Row 33
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.06%  "

# Row 34
$ws.Range("E34").Value = "  -0.36%  "

# Row 35
$ws.Range("D35").Value = "'6.81"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.30%  "

# Row 36
$ws.Range("E36").Value = "  +4.47%  "

# Row 37
$ws.Range("D37").Value = "'162.73"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.02%  "

# Row 38
$ws.Range("E38").Value = "  -2.05%  "

# Row 39
$ws.Range("D39").Value = "'27.20"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.39%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'4.49"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.11%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'2.757.46"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.49%  "

# Row 43
$ws.Range("D43").Value = "'6.24"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.58%  "

# Row 44
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "'0.0672"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.14%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'40.17"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.70%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'24.80"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.98%  "

# Row 47
$ws.Range("E47").Value = "  -1.61%  "

# Row 48
$ws.Range("D48").Value = "'327.58"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.25%  "

# Row 49
$ws.Range("E49").Value = "  -0.55%  "

# Row 50
$ws.Range("D50").Value = "'0.990"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.56%  "

# Row 51
$ws.Range("E51").Value = "  -1.50%  "

